$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.027.36"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "'2.346.95"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'544.61"
$ws.Range("E5").Value = "  +6.06%  "
$ws.Range("D6").Value = "'134.74"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "'2.344.98"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'5.41"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("E13").Value = "  +6.83%  "
$ws.Range("D14").Value = "'2.762.81"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "'23.56"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "'57.989.77"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "'2.349.50"
$ws.Range("E18").Value = "  +4.11%  "
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "'333.77"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'61.65"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +7.55%  "
$ws.Range("E29").Value = "  +4.84%  "
$ws.Range("D30").Value = "'170.26"
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").Value = "'0.0₃0731"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "'6.13"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("E33").Value = "  +17.42%  "
$ws.Range("D34").Value = "'18.46"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'4.18"
$ws.Range("E37").Value = "  +6.32%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").Value = "'1.64"
$ws.Range("E39").Value = "  +4.33%  "
$ws.Range("D40").Value = "'39.30"
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("D41").Value = "'149.00"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'285.89"
$ws.Range("E43").Value = "  +3.49%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'3.60"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "'19.22"
$ws.Range("E45").Value = "  +5.03%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").Value = "'0.561"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").Value = "'0.382"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'17.59"
$ws.Range("E51").Value = "  +3.07%  "
